# Add a new "Vampire Bat" spell bullet to the GDD, right before the
# existing "We will add more spell types as we think of them later" bullet.
#
# The target bullet paragraph currently starts with a
# <w:lastRenderedPageBreak/> marker immediately followed by the text run.
# In the edited document that marker (and its run) become the start of the
# new "Vampire Bat" paragraph, while the "We will add..." text moves into a
# fresh sibling paragraph (losing the page-break marker). We reconstruct
# both paragraphs explicitly via Range.InsertXML so the split happens at
# exactly the right point instead of relying on Find/Replace, which would
# either leave the marker behind or swallow it.

$d = $word.ActiveDocument

$rng = $d.Range(0, 0)
$found = $rng.Find.Execute(
    "We will add more spell types as we think of them later",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'We will add more spell types...' paragraph"
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$pPr = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>' +
       '<w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/>' +
       '<w:adjustRightInd w:val="0"/>' +
       '<w:spacing w:line="252" w:lineRule="auto"/>' +
       '<w:ind w:left="1080" w:hanging="360"/>' +
       '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' +
       '</w:pPr>'

$runRPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/></w:rPr>'

$vampireBatPara = "<w:p $wNs>$pPr<w:r>$runRPr<w:lastRenderedPageBreak/>" +
    "<w:t>Vampire Bat (Average fire rate, Low damage, fires a projectile that has a low chance to restore health to Wand Man upon hitting an enemy)</w:t>" +
    "</w:r></w:p>"

$weWillAddPara = "<w:p $wNs>$pPr<w:r>$runRPr<w:t>We will add more spell types as we think of th</w:t></w:r>" +
    "<w:r>$runRPr<w:t>em later</w:t></w:r></w:p>"

$rng.InsertXML($vampireBatPara + $weWillAddPara)
